# Auto-generated edit script applying the diff to Pandaemonium_Profits workbook
# Updates computed market-price / profit columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 17544684
$ws.Range("I19").Value = 55555900
$ws.Range("J19").Value = 1045.1538
$ws.Range("K19").Value = 55555900
$ws.Range("L19").Value = 1045.1538
$ws.Range("M19").Value = -55555725
$ws.Range("N19").Value = -1395.1538

# Row 28
$ws.Range("H28").Value = 1022.375
$ws.Range("I28").Value = 424.5
$ws.Range("J28").Value = 1620.25
$ws.Range("K28").Value = 424.5
$ws.Range("L28").Value = 1620.25
$ws.Range("M28").Value = 60.5
$ws.Range("N28").Value = -2590.25

# Row 53
$ws.Range("H53").Value = 155.81818
$ws.Range("I53").Value = 103.125
$ws.Range("J53").Value = 185.92857
$ws.Range("K53").Value = 103.125
$ws.Range("L53").Value = 185.92857
$ws.Range("M53").Value = 533.875
$ws.Range("N53").Value = -1459.92857

# Row 62
$ws.Range("H62").Value = 2417.3333
$ws.Range("I62").Value = 2300.8
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2300.8
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1676.8
$ws.Range("N62").Value = -4248

# Row 65
$ws.Range("H65").Value = 2417.3333
$ws.Range("I65").Value = 2300.8
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 11504
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -8384
$ws.Range("N65").Value = -21240

# Row 125
$ws.Range("H125").Value = 6618.1763
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 6618.1763
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 59563.5867
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -64483.5867

# Row 131
$ws.Range("H131").Value = 3916.5
$ws.Range("I131").Value = 527.7273
$ws.Range("J131").Value = 6783.923
$ws.Range("K131").Value = 1583.1819
$ws.Range("L131").Value = 20351.769
$ws.Range("M131").Value = 3456.8181
$ws.Range("N131").Value = -30431.769

# Row 141
$ws.Range("H141").Value = 2608.3704
$ws.Range("I141").Value = 1983.4546
$ws.Range("K141").Value = 5950.3638
$ws.Range("M141").Value = -770.3638000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1820.6
$ws.Range("I45").Value = 1799.8889
$ws.Range("J45").Value = 2007
$ws.Range("K45").Value = 1799.8889
$ws.Range("L45").Value = 2007
$ws.Range("M45").Value = -1422.8889
$ws.Range("N45").Value = -2761

# Row 63
$ws.Range("H63").Value = 2001.6666
$ws.Range("I63").Value = 2001.6666
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2001.6666
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1315.6666
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 2001.6666
$ws.Range("I66").Value = 2001.6666
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10008.333
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6576.333000000001
$ws.Range("N66").ClearContents()

# Row 74
$ws.Range("H74").Value = 4280.59
$ws.Range("I74").Value = 1879.7307
$ws.Range("J74").Value = 9082.308000000001
$ws.Range("K74").Value = 1879.7307
$ws.Range("L74").Value = 9082.308000000001
$ws.Range("M74").Value = -1005.7307
$ws.Range("N74").Value = -10830.308

# Row 77
$ws.Range("H77").Value = 4280.59
$ws.Range("I77").Value = 1879.7307
$ws.Range("J77").Value = 9082.308000000001
$ws.Range("K77").Value = 9398.6535
$ws.Range("L77").Value = 45411.54000000001
$ws.Range("M77").Value = -5030.6535
$ws.Range("N77").Value = -54147.54000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 24
$ws.Range("H24").Value = 3016
$ws.Range("I24").Value = 3016
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 3016
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -2781
$ws.Range("N24").ClearContents()

# Row 80
$ws.Range("H80").Value = 197.90475
$ws.Range("I80").Value = 100
$ws.Range("J80").Value = 202.8
$ws.Range("K80").Value = 100
$ws.Range("L80").Value = 202.8
$ws.Range("M80").Value = 898
$ws.Range("N80").Value = -2198.8

# Row 83
$ws.Range("H83").Value = 197.90475
$ws.Range("I83").Value = 100
$ws.Range("J83").Value = 202.8
$ws.Range("K83").Value = 500
$ws.Range("L83").Value = 1014
$ws.Range("M83").Value = 4492
$ws.Range("N83").Value = -10998

# Row 105
$ws.Range("H105").Value = 870301.4399999999
$ws.Range("I105").Value = 1203036.6
$ws.Range("J105").Value = 5190
$ws.Range("K105").Value = 1203036.6
$ws.Range("L105").Value = 5190
$ws.Range("M105").Value = -1201289.6
$ws.Range("N105").Value = -8684

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1678.6154
$ws.Range("I16").Value = 1204.4
$ws.Range("J16").Value = 1975
$ws.Range("K16").Value = 1204.4
$ws.Range("L16").Value = 1975
$ws.Range("M16").Value = -917.4000000000001
$ws.Range("N16").Value = -2549

# Row 113
$ws.Range("H113").Value = 1678.6154
$ws.Range("I113").Value = 1204.4
$ws.Range("J113").Value = 1975
$ws.Range("K113").Value = 1204.4
$ws.Range("L113").Value = 1975
$ws.Range("M113").Value = 965.5999999999999
$ws.Range("N113").Value = -6315

# Row 141
$ws.Range("H141").Value = 20296
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 717.64703
$ws.Range("I113").Value = 712.8205
$ws.Range("J113").Value = 733.3333
$ws.Range("K113").Value = 2138.4615
$ws.Range("L113").Value = 2199.9999
$ws.Range("M113").Value = 31.53849999999966
$ws.Range("N113").Value = -6539.9999

# Row 131
$ws.Range("H131").Value = 54806.723
$ws.Range("I131").Value = 2680
$ws.Range("J131").Value = 87978.27
$ws.Range("K131").Value = 8040
$ws.Range("L131").Value = 263934.81
$ws.Range("M131").Value = -3000
$ws.Range("N131").Value = -274014.81

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2914.818
$ws.Range("I126").Value = 1939
$ws.Range("J126").Value = 3472.4285
$ws.Range("K126").Value = 5817
$ws.Range("L126").Value = 10417.2855
$ws.Range("M126").Value = -3347
$ws.Range("N126").Value = -15357.2855

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 360.08334
$ws.Range("I55").Value = 268.77777
$ws.Range("J55").Value = 634
$ws.Range("K55").Value = 268.77777
$ws.Range("L55").Value = 634
$ws.Range("M55").Value = -95.77776999999998
$ws.Range("N55").Value = -980

# Row 93
$ws.Range("H93").Value = 221.83333
$ws.Range("I93").Value = 176.75
$ws.Range("K93").Value = 176.75
$ws.Range("M93").Value = 1071.25

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 4512.8887
$ws.Range("I113").Value = 8705.166999999999
$ws.Range("J113").Value = 1159.0667
$ws.Range("K113").Value = 26115.501
$ws.Range("L113").Value = 3477.2001
$ws.Range("M113").Value = -23945.501
$ws.Range("N113").Value = -7817.2001

# Row 122
$ws.Range("H122").Value = 2261.8462
$ws.Range("I122").Value = 1764
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 5292
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -2842
$ws.Range("N122").Value = -19900
